$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff:
# A2: 0 -> 1
# B2: 95 -> 90
# A3: 1 -> 0
# B3: 69 -> 74
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 90
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 74
